# "Generate Report for Handoff"
# The localization-status report was regenerated: the zh-cn / de-de entries
# moved from "In Translation" to "Ready for handoff", the handoff timestamps
# were refreshed, and the Status column widened to fit the new (longer)
# status text on all three sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-25 08:42:14"

# Widen the (now longer) status columns to fit "Ready for handoff"
$ws1.Columns.Item(5).ColumnWidth = 16.3
$ws1.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-25 08:42:07"
$ws2.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ------------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-25 08:42:14"
$ws3.Columns.Item(3).ColumnWidth = 16.3
